$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D ("Price") holds plain text in the source data -- several values
# look numeric (e.g. "0.998", "68.30", "0.0420") and Excel would silently
# coerce/round them if assigned directly, so force text format per-cell first.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.416.05'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.921.95'
$ws.Range('E3').Value = '  -2.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '376.93'
$ws.Range('E5').Value = '  +6.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.78'
$ws.Range('E6').Value = '  -4.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.542'
$ws.Range('E7').Value = '  -2.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.584'
$ws.Range('E9').Value = '  -4.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.96'
$ws.Range('E10').Value = '  -3.14%  '
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0834'
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('E13').Value = '  -3.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.379.06'
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.33'
$ws.Range('E15').Value = '  -3.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.916.83'
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.926'
$ws.Range('E17').Value = '  -8.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.279.26'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.34'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.93'
$ws.Range('E21').Value = '  -4.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0945'
$ws.Range('E22').Value = '  -2.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.30'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.82'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E26').Value = '  -4.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.12'
$ws.Range('E27').Value = '  -4.87%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.24'
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.85'
$ws.Range('E31').Value = '  +7.44%  '
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.79'
$ws.Range('E33').Value = '  -4.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.19'
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '33.97'
$ws.Range('E36').Value = '  -5.75%  '
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0420'
$ws.Range('E38').Value = '  -4.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.99'
$ws.Range('E39').Value = '  -10.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.91'
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.53'
$ws.Range('E41').Value = '  -11.16%  '
$ws.Range('E42').Value = '  -7.94%  '
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '122.31'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.47'
$ws.Range('E45').Value = '  -6.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.06'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.269'
$ws.Range('E47').Value = '  +10.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.025.56'
$ws.Range('E48').Value = '  -4.57%  '
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('E50').Value = '  -5.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.202.99'
$ws.Range('E51').Value = '  -2.85%  '
